# Replace the two-line sample text ("some text" / "hello qa.guru students!")
# with a single line of text ("test line") in cell B4, and remove the old
# second line that lived in B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 keeps its existing style (s="1") but gets new text.
$ws.Range("B4").Value = "test line"

# B5 no longer holds any data - clear it out entirely.
$ws.Range("B5").ClearContents()

# Leave the active cell on B5, matching the saved selection state.
$ws.Range("B5").Select()
